$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Mean (H) and Std (I) columns for PreferenceOrder prediction-type rows
$values = @{
    26 = @{ H = 0.65627; I = 0.03294 }
    27 = @{ H = 0.24075; I = 0.06852 }
    28 = @{ H = 0.66917; I = 0.03294 }
    29 = @{ H = 0.30818; I = 0.07199 }
    30 = @{ H = 0.65678; I = 0.03279 }
    31 = @{ H = 0.23998; I = 0.06542000000000001 }
    32 = @{ H = 0.67103; I = 0.0332 }
    33 = @{ H = 0.3101; I = 0.07198 }
    34 = @{ H = 0.56556; I = 0.01855 }
    35 = @{ H = 0.03121; I = 0.02742 }
    36 = @{ H = 0.56555; I = 0.01833 }
    37 = @{ H = 0.03314; I = 0.02684 }
    38 = @{ H = 0.56556; I = 0.01855 }
    39 = @{ H = 0.03121; I = 0.02742 }
    40 = @{ H = 0.56555; I = 0.01833 }
    41 = @{ H = 0.03314; I = 0.02684 }
    66 = @{ H = 0.55416; I = 0.02295 }
    67 = @{ H = 0.10214; I = 0.02493 }
    68 = @{ H = 0.58296; I = 0.02323 }
    69 = @{ H = 0.1915; I = 0.03782 }
    70 = @{ H = 0.5563399999999999; I = 0.02176 }
    71 = @{ H = 0.10175; I = 0.02495 }
    72 = @{ H = 0.59002; I = 0.02511 }
    73 = @{ H = 0.20575; I = 0.03822 }
    74 = @{ H = 0.54965; I = 0.01838 }
    75 = @{ H = 0.03892; I = 0.01972 }
    76 = @{ H = 0.55132; I = 0.02049 }
    77 = @{ H = 0.05241; I = 0.02921 }
    78 = @{ H = 0.55016; I = 0.0182 }
    79 = @{ H = 0.03892; I = 0.01972 }
    80 = @{ H = 0.5510699999999999; I = 0.0203 }
    81 = @{ H = 0.05203; I = 0.03008 }
    106 = @{ H = 0.63123; I = 0.03341 }
    107 = @{ H = 0.19689; I = 0.07252 }
    108 = @{ H = 0.6442099999999999; I = 0.0323 }
    109 = @{ H = 0.25665; I = 0.07267999999999999 }
    110 = @{ H = 0.63123; I = 0.03417 }
    111 = @{ H = 0.19535; I = 0.07038 }
    112 = @{ H = 0.6462; I = 0.03155 }
    113 = @{ H = 0.26321; I = 0.07248 }
    114 = @{ H = 0.5524; I = 0.01924 }
    115 = @{ H = 0.0235; I = 0.02508 }
    116 = @{ H = 0.55227; I = 0.01959 }
    117 = @{ H = 0.02388; I = 0.02581 }
    118 = @{ H = 0.5523400000000001; I = 0.0194 }
    119 = @{ H = 0.0235; I = 0.02508 }
    120 = @{ H = 0.55208; I = 0.01948 }
    121 = @{ H = 0.02388; I = 0.02581 }
    146 = @{ H = 0.58479; I = 0.03008 }
    147 = @{ H = 0.1476; I = 0.05233 }
    148 = @{ H = 0.60431; I = 0.03027 }
    149 = @{ H = 0.2301; I = 0.061 }
    150 = @{ H = 0.58512; I = 0.02927 }
    151 = @{ H = 0.14722; I = 0.0519 }
    152 = @{ H = 0.61099; I = 0.02639 }
    153 = @{ H = 0.24281; I = 0.05834 }
    154 = @{ H = 0.55768; I = 0.02735 }
    155 = @{ H = 0.04971; I = 0.03027 }
    156 = @{ H = 0.5598; I = 0.02884 }
    157 = @{ H = 0.06281; I = 0.03314 }
    158 = @{ H = 0.55781; I = 0.02725 }
    159 = @{ H = 0.04971; I = 0.03027 }
    160 = @{ H = 0.55974; I = 0.0285 }
    161 = @{ H = 0.06319; I = 0.03291 }
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item([int]$row, 8).Value = $values[$row].H
    $ws.Cells.Item([int]$row, 9).Value = $values[$row].I
}
